# Append the 2025-04-25 price row to each of the 9 price sheets in the
# workbook. Each sheet currently ends at row 54 (dimension A1:B54); a new
# row 55 is added carrying the same price as row 54 (the most recent
# Argent/solar price), with the date moved forward one day.

$wb = $excel.ActiveWorkbook

# Map of sheet name -> new price value for column B, row 55.
# Column A (date) is the same for every sheet: 2025-04-25.
$newRow = 55
$newDate = "2025-04-25"

$prices = @{
    "N-Dense"                    = "39.5"
    "N-Type"                     = "39.31"
    "N-type Wafer"               = "1.18"
    "Cell Topcon 183mm"          = "0.288"
    "Module Topcon 183mm"        = "0.09"
    "Silver Rear_side"           = "5,419"
    "Silver Busbar front-side"   = "8,113"
    "Silver finger front-side"   = "8,163"
    "USD_CNY"                    = "7.3088"
}

foreach ($ws in $wb.Worksheets) {
    $name = $ws.Name
    if ($prices.ContainsKey($name)) {
        $price = $prices[$name]

        $dateCell = $ws.Cells.Item($newRow, 1)
        $dateCell.NumberFormat = "@"
        $dateCell.Value = $newDate

        $priceCell = $ws.Cells.Item($newRow, 2)
        $priceCell.NumberFormat = "@"
        $priceCell.Value = $price
    }
}
